{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// The three trailing paragraphs that need to be removed:\n//   1) an empty paragraph\n//   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   3) \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n// They immediately follow the \"Textos fornecidos...\" paragraph and are\n// immediately followed by another empty paragraph (kept) + a page-break paragraph (kept).\nconst items = paragraphs.items;\n\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"Textos fornecidos pelo professor da disciplina\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex !== -1) {\n  const toDelete = [];\n  // empty paragraph right after the anchor\n  if (items[anchorIndex + 1] && items[anchorIndex + 1].text === \"\") {\n    toDelete.push(items[anchorIndex + 1]);\n  }\n  // \"Ver no Jupiter...\" paragraph\n  if (items[anchorIndex + 2] && items[anchorIndex + 2].text.indexOf(\"Ver no Jupiter\") !== -1) {\n    toDelete.push(items[anchorIndex + 2]);\n  }\n  // \"\u00a9 2020 ...\" paragraph\n  if (items[anchorIndex + 3] && items[anchorIndex + 3].text.indexOf(\"Creative Commons Attribution\") !== -1) {\n    toDelete.push(items[anchorIndex + 3]);\n  }\n\n  for (const p of toDelete) {\n    p.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"Textos fornecidos pelo professor da disciplina ...\" paragraph;\n# the three paragraphs that immediately follow it are the ones to remove:\n#   1) an empty paragraph\n#   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   3) \"(c) 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n#       pages. Original theme under Creative Commons Attribution\"\n# A further empty paragraph and a page-break paragraph follow afterwards and\n# must be left untouched.\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*Textos fornecidos pelo professor da disciplina*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -ge 1) {\n    # Walk the three target paragraphs from last to first so deleting one\n    # doesn't shift the index of the ones still pending.\n    for ($offset = 3; $offset -ge 1; $offset--) {\n        $idx = $anchorIndex + $offset\n        if ($idx -le $d.Paragraphs.Count) {\n            $d.Paragraphs.Item($idx).Range.Delete()\n        }\n    }\n}\n"}
